# Regenerate orders with updated distance/sizes
# Applies value substitutions across the whole used range:
#   D51 -> D55,  D64 -> D69,  D80 -> D86   (Distance codes)
#   S30 -> S31                              (Size code; S20/S25 unchanged)
# The substitutions are applied to the token as a whole (not inside a
# longer number), mirroring how these tokens appear embedded in strings
# like "Face18_D51_S30", "Face18_D51_S30_l.png", "Fixation_D51_l.png", "D51".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
  for ($c = 1; $c -le $colCount; $c++) {
    $cell = $ws.Cells.Item($r, $c)
    $val = $cell.Value2

    if ($val -is [string]) {
      $newVal = $val
      $newVal = $newVal -replace 'D51(?!\d)', 'D55'
      $newVal = $newVal -replace 'D64(?!\d)', 'D69'
      $newVal = $newVal -replace 'D80(?!\d)', 'D86'
      $newVal = $newVal -replace 'S30(?!\d)', 'S31'

      if ($newVal -ne $val) {
        $cell.Value2 = $newVal
      }
    }
  }
}
